# Applies the content changes described in the commit "Update dagstuhl
# talk and pc memberships."
#
# The underlying XML diff also wraps many existing (unchanged) phrases in
# <w:proofErr> start/end markers -- that is Word's grammar-checker churn
# (from opening the file in a newer Word build) and carries no visible or
# semantic content change, so it is not reproduced here; only actual text
# changes are applied.

$d = $word.ActiveDocument

# 1. Invited talk: Dagstuhl workshop entry was rewritten (new talk title,
#    dropped the workshop name, and the day changed from the 4th to the 6th).
$d.Content.Find.Execute(
    "Approximating Betweenness Centrality through Sampling with the Rademacher Averages, Dagstuhl Workshop on Probabilistic Methods in the Design and Analysis of Algorithms, Schloss Dagstuhl (Germany), April 4",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Rademacher Averages: Theory and Practice, Schloss Dagstuhl (Germany), April 6",
    2) | Out-Null

# 2. Program Committees: add ACM CIKM '17 membership.
$d.Content.Find.Execute(
    "ACM CIKM ‘16, ‘15, ‘14",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "ACM CIKM ‘17, ‘16, ‘15, ‘14",
    2) | Out-Null

# 3. Program Committees: add IEEE ICDE '18 membership.
$d.Content.Find.Execute(
    "IEEE ICDE ‘17",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "IEEE ICDE ‘18, ‘17",
    2) | Out-Null

# 4. Footer page-count field: cached result bumped from 4 to 7 (document
#    grew by a few lines/paragraphs after the above proofing pass).
$sec = $d.Sections(1)
$ftr = $sec.Footers(1)
foreach ($f in $ftr.Range.Fields) {
    if ($f.Code.Text.Trim() -eq "PAGE") {
        $f.Result.Text = "7"
    }
}
